$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '41.264.00'
Set-TextValue 'E2' '  -3.75%  '

Set-TextValue 'D3' '2.462.42'
Set-TextValue 'E3' '  -2.91%  '

Set-TextValue 'E4' '  +0.05%  '

Set-TextValue 'D5' '311.50'
Set-TextValue 'E5' '  -0.09%  '

Set-TextValue 'D6' '94.16'
Set-TextValue 'E6' '  -6.55%  '

Set-TextValue 'D7' '0.548'
Set-TextValue 'E7' '  -3.53%  '

Set-TextValue 'D9' '0.498'
Set-TextValue 'E9' '  -4.99%  '

Set-TextValue 'D10' '33.39'
Set-TextValue 'E10' '  -6.75%  '

Set-TextValue 'D11' '0.0780'
Set-TextValue 'E11' '  -3.27%  '

Set-TextValue 'D12' '0.109'

Set-TextValue 'D13' '6.98'
Set-TextValue 'E13' '  -4.82%  '

Set-TextValue 'D14' '2.842.78'
Set-TextValue 'E14' '  -2.88%  '

Set-TextValue 'D15' '2.474.11'
Set-TextValue 'E15' '  -2.53%  '

Set-TextValue 'D16' '14.98'
Set-TextValue 'E16' '  -2.74%  '

Set-TextValue 'D17' '0.786'
Set-TextValue 'E17' '  -3.99%  '

Set-TextValue 'D18' '41.272.40'
Set-TextValue 'E18' '  -3.77%  '

Set-TextValue 'D19' '6.30'
Set-TextValue 'E19' '  -5.64%  '

Set-TextValue 'D20' '0.0₃0921'
Set-TextValue 'E20' '  -3.53%  '

Set-TextValue 'D21' '11.18'
Set-TextValue 'E21' '  -9.47%  '

Set-TextValue 'D22' '68.55'
Set-TextValue 'E22' '  -2.04%  '

Set-TextValue 'D23' '236.89'
Set-TextValue 'E23' '  -2.94%  '

Set-TextValue 'E24' '  -4.78%  '

Set-TextValue 'E25' '  +0.13%  '

Set-TextValue 'D26' '1.90'
Set-TextValue 'E26' '  -6.90%  '

Set-TextValue 'D27' '24.00'
Set-TextValue 'E27' '  -6.07%  '

Set-TextValue 'E28' '  -4.84%  '

Set-TextValue 'D29' '9.64'
Set-TextValue 'E29' '  -5.60%  '

Set-TextValue 'D30' '36.33'
Set-TextValue 'E30' '  -6.13%  '

Set-TextValue 'D31' '151.54'
Set-TextValue 'E31' '  -4.70%  '

Set-TextValue 'D32' '5.48'
Set-TextValue 'E32' '  -6.37%  '

Set-TextValue 'D33' '2.63'
Set-TextValue 'E33' '  -6.05%  '

Set-TextValue 'D34' '2.59'
Set-TextValue 'E34' '  -3.40%  '

Set-TextValue 'D35' '0.0745'
Set-TextValue 'E35' '  -6.21%  '

Set-TextValue 'E36' '  -3.22%  '

Set-TextValue 'D37' '17.10'
Set-TextValue 'E37' '  -6.89%  '

Set-TextValue 'E38' '  -5.53%  '

Set-TextValue 'E39' '  -3.24%  '

Set-TextValue 'D40' '0.102'
Set-TextValue 'E40' '  -8.28%  '

Set-TextValue 'D41' '4.22'
Set-TextValue 'E41' '  +0.79%  '

Set-TextValue 'E42' '  +0.15%  '

Set-TextValue 'D43' '19.58'
Set-TextValue 'E43' '  -10.33%  '

Set-TextValue 'D44' '1.985.79'
Set-TextValue 'E44' '  -0.72%  '

Set-TextValue 'D45' '0.0284'
Set-TextValue 'E45' '  -5.29%  '

Set-TextValue 'D46' '3.01'
Set-TextValue 'E46' '  -9.41%  '

Set-TextValue 'E47' '  -4.41%  '

Set-TextValue 'D48' '2.708.98'
Set-TextValue 'E48' '  -2.54%  '

Set-TextValue 'D49' '69.19'
Set-TextValue 'E49' '  -4.62%  '

Set-TextValue 'D50' '96.36'
Set-TextValue 'E50' '  -5.02%  '

Set-TextValue 'D51' '74.41'
Set-TextValue 'E51' '  -7.14%  '
